$d = $word.ActiveDocument

# Locate the end of the "LOQ4073: Química Geral II (Requisito fraco)" paragraph
# (i.e. the position right after its paragraph mark / start of the next paragraph).
$startRng = $d.Content
$startRng.Find.Execute("LOQ4073: Química Geral II (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startRng.Expand(4) | Out-Null
$startPos = $startRng.End

# Locate the end of the copyright/footer paragraph
# ("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ...")
# which, together with the blank paragraph and the "Ver no Jupiter ..." paragraph
# right after LOQ4073, must be removed.
$endRng = $d.Content
$endRng.Find.Execute("Powered by Jekyll and Github pages", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endRng.Expand(4) | Out-Null
$endPos = $endRng.End

# Delete the three trailing paragraphs (blank line, "Ver no Jupiter ..." line,
# and the "© 2020 ..." line) in one shot, leaving the LOQ4073 paragraph followed
# directly by the original trailing blank paragraph / page break.
$d.Range($startPos, $endPos).Delete()
